$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Aufgabe "Methode zum Anzeigen von Verfuegbbaren Flugzeugen"):
# Bearbeiter-Feld aktualisieren - Leon hat den Fix gemacht.
$ws.Range("B5").Value = "Leroy | Fixed by Leon"

# Row 9 (Aufgabe "Methode zum hinzufuegen einer leihe."):
# Status von "Start" (rot) auf "Fertig" (gruen) setzen - Format von einer
# bereits fertigen Zeile (C2) uebernehmen, dann den Text anpassen.
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "Fertig"

# Aktuelle Auswahl wie im Original-Commit auf B15 setzen.
$null = $ws.Range("B15").Select()
